$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The rows of the player table (A2:C17) were reordered (e.g. by dragging /
# cutting-and-pasting rows around). Re-write the whole table body in its
# new row order so every (player, position, team) record ends up on the
# correct row.

$data = @(
    @("Ja Morant",           "PG",       "Memphis Grizzlies"),
    @("Luka Doncic",         "PG,SG",    "Dallas Mavericks"),
    @("De'Aaron Fox",        "PG",       "Sacramento Kings"),
    @("Vasilije Micic",      "PG,SG",    "Charlotte Hornets"),
    @("DeMar DeRozan",       "SF,PF",    "Sacramento Kings"),
    @("Evan Mobley",         "PF,C",     "Cleveland Cavaliers"),
    @("Amen Thompson",       "SG,SF",    "Houston Rockets"),
    @("Bennedict Mathurin",  "SG,SF",    "Indiana Pacers"),
    @("Santi Aldama",        "PF,C",     "Memphis Grizzlies"),
    @("Brook Lopez",         "C",        "Milwaukee Bucks"),
    @("Mikal Bridges",       "SG,SF,PF", "New York Knicks"),
    @("Scottie Barnes",      "SG,SF,PF", "Toronto Raptors"),
    @("Tyler Herro",         "PG,SG",    "Miami Heat"),
    @("Josh Giddey",         "PG,SG,SF", "Chicago Bulls"),
    @("Nikola Vucevic",      "PF,C",     "Chicago Bulls"),
    @("Miles Bridges",       "SF,PF",    "Charlotte Hornets")
)

$row = 2
foreach ($rec in $data) {
    $ws.Range("A$row").Value = $rec[0]
    $ws.Range("B$row").Value = $rec[1]
    $ws.Range("C$row").Value = $rec[2]
    $row++
}
